$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1 - copy style from existing header (e.g. E1) then set values
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Match the styling of the existing header row (bold, centered, bordered)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Restore the header text after paste (paste formats only touches formatting,
# but ensure values remain correct)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Boolean data for rows 2-12 across columns F, G, H - all FALSE except G6 = TRUE
$ws.Range("F2:H12").Value = $false
$ws.Range("G6").Value = $true
